$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new price-tracking row (row 38) with:
#   A38 = "2026-02-07" (new date, new shared string)
#   B38 = "100000"     (price, reuses the existing shared string from B37)
#   C38 = "0"          (discount, reuses the existing shared string from C37)
#   D38 = "0"          (incredible flag, reuses the existing shared string from D37)
#
# Plain `Range.Value = "..."` assignment lets Excel "smart type" the input
# (dates/numbers get coerced to real numeric/date cells with a new
# NumberFormat style), which doesn't match how this sheet stores its data
# (everything as plain text / shared strings with no special styling).
# Instead we stage the literal text in a scratch cell as a text FORMULA
# (guaranteed text, never re-interpreted) and paste just the resulting
# VALUE into place, which keeps it as plain text with no style change.

$scratch = $ws.Range("ZZ1")
$scratch.Formula = '="2026-02-07"'
$scratch.Copy()
$ws.Range("A38").PasteSpecial(-4163)
$scratch.Clear()

# B38/C38/D38 already have identical text values earlier in the column, so
# just copy those existing text cells straight down - this keeps them as
# shared-string text cells without touching any number formatting.
$ws.Range("B37").Copy($ws.Range("B38"))
$ws.Range("C37").Copy($ws.Range("C38"))
$ws.Range("D37").Copy($ws.Range("D38"))

Write-Output "Added row 38: $($ws.Range('A38').Value2), $($ws.Range('B38').Value2), $($ws.Range('C38').Value2), $($ws.Range('D38').Value2)"
